# Update the AasCode regular-expression patterns in the "Conventies" sheet
# for Lift / Roltrap / Rolpad / Hellingbaan: the numeric suffix changes from
# an exact 3-digit group ("\d{3}") to a 1-to-4-digit group ("\d{1,4}").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conventies")

$ws.Range("D5").Value  = "^L\d{1,4}$"
$ws.Range("D6").Value  = "^RT\d{1,4}$"
$ws.Range("D7").Value  = "^RP\d{1,4}$"
$ws.Range("D8").Value  = "^HB\d{1,4}$"
$ws.Range("D10").Value = "^RT\d{1,4}$"

# Widen column I (Assetbeschrijving - Voorbeeld) to fit its content.
$ws.Columns.Item(9).ColumnWidth = 70

# Move the frozen-pane viewport / selection on the bottom-right pane back to
# the top of the table (D13), which also resets the scrolled top-left cell.
$ws.Activate()
$ws.Range("D13").Select()
